$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BookATicket")

# Add the new "date" / "current date" column (H) to the BookATicket sheet
$ws.Range("H1").Value = "date"
$ws.Range("H2").Value = "current date"

# Give the new column a sensible width (auto-fit to its contents)
$ws.Columns("H").ColumnWidth = 11.85546875

# Style H1 like the other header cells (bold font + yellow fill) but with a
# left/right only border, then extend the same border down to H2.
$ws.Range("H1").Borders.Item(7).LineStyle = 1
$ws.Range("H1").Borders.Item(10).LineStyle = 1
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").Interior.Color = 65535

$ws.Range("H2").Borders.Item(7).LineStyle = 1
$ws.Range("H2").Borders.Item(10).LineStyle = 1

# Update the selection to match the edited workbook
$null = $ws.Range("K6").Select()
